# Manchies Registration Page - planning notes update
#
# Adds two new paragraphs after the existing "Registration dialog - will
# use the exact same design." paragraph:
#   1. A blank paragraph (same tab-stop / run-formatting as the rest of
#      the body text).
#   2. A paragraph with the text "This page general design is going to
#      be the same as the login page, but with other fields." - this
#      paragraph also inherits the "_GoBack" bookmark that used to sit
#      at the end of the previous (now no longer last) paragraph, since
#      the bookmark always tracks the most recent edit location.

$d = $word.ActiveDocument

# The paragraph that currently ends the document body ("... use the
# exact same design.") - grab it before we start inserting anything.
$lastPara = $d.Paragraphs.Last

# --- 1. Blank paragraph -----------------------------------------------
$lastPara.Range.InsertParagraphAfter()

# --- 2. New paragraph with the planning note --------------------------
$midPara = $d.Paragraphs.Last
$midPara.Range.InsertParagraphAfter()

$finalPara = $d.Paragraphs.Last
# Insert the sentence plus a throw-away trailing marker character; the
# marker lets us park a zero-width bookmark range that is NOT sitting at
# the very end of the document/paragraph (which this host snaps/expands
# to cover a whole run) - we delete the marker once the bookmark has
# been anchored.
$finalPara.Range.InsertBefore("This page general design is going to be the same as the login page, but with other fields.#")

# --- 3. Relocate the "_GoBack" bookmark to the new last paragraph -----
$finalPara = $d.Paragraphs.Last
$markerRange = $finalPara.Range.Duplicate
$markerRange.MoveEnd(1, -1)                 # drop the paragraph mark
$markerRange.MoveStart(1, $markerRange.Text.Length - 1)  # just the "#"
$markerRange.Collapse(1)                    # collapse to before "#"
$d.Bookmarks.Add("_GoBack", $markerRange)

# --- 4. Remove the temporary marker character --------------------------
$finalPara = $d.Paragraphs.Last
$markerRange2 = $finalPara.Range.Duplicate
$markerRange2.MoveEnd(1, -1)
$markerRange2.MoveStart(1, $markerRange2.Text.Length - 1)
$markerRange2.Delete()
